$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(214).Insert()

$ws.Cells.Item(214, 1).Value = 7
$ws.Cells.Item(214, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(214, 3).Value = "Ñuble"
$ws.Cells.Item(214, 4).Value = 44943
$ws.Cells.Item(214, 5).Value = 16
$ws.Cells.Item(214, 6).Value = 100112024
$ws.Cells.Item(214, 7).Value = "Choclo"
$ws.Cells.Item(214, 8).Value = "Choclero"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 30000
$ws.Cells.Item(214, 11).Value = 250
$ws.Cells.Item(214, 12).Value = 300
$ws.Cells.Item(214, 13).Value = 275
$ws.Cells.Item(214, 14).Value = "`$/unidad"
$ws.Cells.Item(214, 15).Value = "Región del Maule"
$ws.Cells.Item(214, 16).Value = 275
$ws.Cells.Item(214, 17).Value = 1
$ws.Cells.Item(214, 18).Value = "Hortaliza"
